# Fruta / hortaliza, semanal
# Insert two new price records (rows 281-282) for "Navel Late" oranges
# dated 44504, pushing the existing rows 281-333 down to 283-335.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 281.
$ws.Rows.Item(281).Insert()
$ws.Rows.Item(281).Insert()

# New row 281
$ws.Cells.Item(281, 1).Value = 5
$ws.Cells.Item(281, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(281, 3).Value = "Maule"
$ws.Cells.Item(281, 4).Value = 44504
$ws.Cells.Item(281, 5).Value = 7
$ws.Cells.Item(281, 6).Value = "Fruta"
$ws.Cells.Item(281, 7).Value = 100102
$ws.Cells.Item(281, 8).Value = "Cítricos"
$ws.Cells.Item(281, 9).Value = 100102005
$ws.Cells.Item(281, 10).Value = "Naranja"
$ws.Cells.Item(281, 11).Value = "Navel Late"
$ws.Cells.Item(281, 12).Value = "Primera"
$ws.Cells.Item(281, 13).Value = 300
$ws.Cells.Item(281, 14).Value = 8000
$ws.Cells.Item(281, 15).Value = 8000
$ws.Cells.Item(281, 16).Value = 8000
$ws.Cells.Item(281, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(281, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(281, 19).Value = 533
$ws.Cells.Item(281, 20).Value = 15

# New row 282
$ws.Cells.Item(282, 1).Value = 5
$ws.Cells.Item(282, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(282, 3).Value = "Maule"
$ws.Cells.Item(282, 4).Value = 44504
$ws.Cells.Item(282, 5).Value = 7
$ws.Cells.Item(282, 6).Value = "Fruta"
$ws.Cells.Item(282, 7).Value = 100102
$ws.Cells.Item(282, 8).Value = "Cítricos"
$ws.Cells.Item(282, 9).Value = 100102005
$ws.Cells.Item(282, 10).Value = "Naranja"
$ws.Cells.Item(282, 11).Value = "Navel Late"
$ws.Cells.Item(282, 12).Value = "Primera"
$ws.Cells.Item(282, 13).Value = 550
$ws.Cells.Item(282, 14).Value = 7500
$ws.Cells.Item(282, 15).Value = 8000
$ws.Cells.Item(282, 16).Value = 7682
$ws.Cells.Item(282, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(282, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(282, 19).Value = 512
$ws.Cells.Item(282, 20).Value = 15
